# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap two pairs of country names that moved position in the shared-string
# table (the worksheet cells keep pointing at the same row/index, so the
# visible text at those two rows swaps) ---
$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("A201").Value = "Belice"

$ws.Range("A215").Value = "San Bartolome"
$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"

# --- Update "last updated" timestamp text ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Mayo de 2020 a las 21:10"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1782616
$ws.Range("C4").Value = 14155
$ws.Range("D4").Value = 502263
$ws.Range("E4").Value = 1176355
$ws.Range("G4").Value = 668
$ws.Range("H4").Value = 103998

# --- Row 5: Brasil ---
$ws.Range("B5").Value = 443876
$ws.Range("C5").Value = 5064
$ws.Range("E5").Value = 223794
$ws.Range("G5").Value = 137
$ws.Range("H5").Value = 26901

# --- Row 53: Barein ---
$ws.Range("B53").Value = 10449
$ws.Range("C53").Value = 397
$ws.Range("D53").Value = 5700
$ws.Range("E53").Value = 4734

# --- Row 64: Marruecos ---
$ws.Range("B64").Value = 7714
$ws.Range("C64").Value = 71
$ws.Range("D64").Value = 5271
$ws.Range("E64").Value = 2241

# --- Row 200 / 201 numeric data follows the swapped country names ---
$ws.Range("D200").Value = 18
$ws.Range("H200").Value = 0

$ws.Range("D201").Value = 16
$ws.Range("H201").Value = 2
